$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Sep 08 18:13:00 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 08 18:13:14 EDT 2023"
$ws.Range("D3").Value = "Extension Payments"
$ws.Range("B4").Value = "Fri Sep 08 18:13:29 EDT 2023"

$ws.Range("D3").Select()
